$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.493.90'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.97%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.460.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '90.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.534'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.78%  '

$ws.Range("E8").Value = '  +0.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '31.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0770'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.16%  '

$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.821.25'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.12'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.90%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.415.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.756'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.213.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0906'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.85'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.96%  '

$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '150.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0748'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.77%  '

$ws.Range("E35").Value = '  -3.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.61'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.45%  '

$ws.Range("E39").Value = '  -3.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0993'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.96%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.22%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.927.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0275'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.678.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '94.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.173'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.02%  '
